# ANA1: Added information on the outliers from the 2nd Analysis Phase.
# This covered GeoTools, Findbugs and NetBeans.
#
# The "3rdParty" table is re-sorted by "Total Number of Unique
# Vulnerabilities" (column D) descending, so the outlier projects show up
# at the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3rdParty")

# Sort the data range (header in row 1, data in rows 2-77) by column D
# (Total Number of Unique Vulnerabilities), descending.
$sortRange = $ws.Range("A1:H77")
$sortKey   = $ws.Range("D1:D77")
$sortRange.Sort($sortKey, 2, $null, $null, 1, $null, 1, 1, $false, $null, $null, 1)

# Bring the 3rdParty sheet to the front, scroll back to the top and select
# the first data cell under the sorted column so the view matches the
# freshly-sorted table.
$ws.Activate()
$ws.Range("D2").Select()
